# Apply updated "dSF" (column F) values to rows 2-29 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = 2
    4  = -3
    5  = -3
    7  = -3
    9  = 5
    10 = 2
    11 = -1
    12 = -4
    13 = 2
    14 = 4
    15 = 1
    16 = -2
    17 = 2
    18 = 3
    19 = -4
    20 = -2
    21 = -4
    22 = 8
    24 = 3
    26 = -4
    27 = -3
    28 = -6
    29 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
